$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data cells H20:H34 (AIC values copied from H2:H16) ---
# Build the target style on H20 first: copy number formatting base from G20
# (fontId already matches Times New Roman / theme color used throughout the
# table), then switch the number format to 2 decimals and right-align it.
$ws.Range("G20").Copy()
$ws.Range("H20").PasteSpecial(-4122)
$ws.Range("H20").Value = $ws.Range("H2").Value2
$ws.Range("H20").NumberFormat = "0.00"
$ws.Range("H20").HorizontalAlignment = -4152

# --- Header cell H19 ("AIC") ---
# Start again from G20's format (keeps the correct font), reset the number
# format back to General, add the top/bottom border used by the other
# header cells, and right-align it.
$ws.Range("G20").Copy()
$ws.Range("H19").PasteSpecial(-4122)
$ws.Range("H19").Value = "AIC"
$ws.Range("H19").NumberFormat = "general"
$ws.Range("H19").Borders.Item(9).LineStyle = 1
$ws.Range("H19").Borders.Item(8).LineStyle = 1
$ws.Range("H19").HorizontalAlignment = -4152

# --- Remaining data cells H21:H34 ---
# Reuse H20's already-built style by copying formats from it (avoids
# creating redundant style table entries).
for ($i = 21; $i -le 34; $i++) {
    $srcRow = $i - 18
    $ws.Range("H20").Copy()
    $ws.Range("H" + $i).PasteSpecial(-4122)
    $ws.Range("H" + $i).Value = $ws.Range("H" + $srcRow).Value2
}

# --- Selection moved to I19 ---
$ws.Range("I19").Select()
